$wb = $excel.ActiveWorkbook

# "Orders" sheet (4th sheet / sheet4.xml) gets a new "PredefinedTime" column
$ws = $wb.Worksheets.Item(4)
$ws.Activate()

# Add the new header cell in column I, row 1
$ws.Range("I1").Value = "PredefinedTime"

# Match the header formatting used by the rest of row 1 (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Size the new column to fit its header content
$ws.Columns.Item(9).AutoFit()

# Leave the new header cell selected, matching the saved view state
[void]$ws.Range("I1").Select()
